# Weekly update: insert two new report rows (week of 2022-05-13, serial 44694)
# at the top of the Tomate / Vega Monumental Concepción block, pushing the
# existing rows (old 352..415) down by two (to 354..417).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 352.
$ws.Rows("352:353").Insert()

# Row 352: Tomate, Larga vida, Primera
$ws.Range("A352").Value = 11
$ws.Range("B352").Value = 'Vega Monumental Concepción'
$ws.Range("C352").Value = 'Bíobío'
$ws.Range("D352").Value = 44694
$ws.Range("E352").Value = 8
$ws.Range("F352").Value = 100112020
$ws.Range("G352").Value = 'Tomate'
$ws.Range("H352").Value = 'Larga vida'
$ws.Range("I352").Value = 'Primera'
$ws.Range("J352").Value = 800
$ws.Range("K352").Value = 19000
$ws.Range("L352").Value = 20000
$ws.Range("M352").Value = 19500
$ws.Range("N352").Value = '$/bandeja 18 kilos'
$ws.Range("O352").Value = 'Región de Arica y Parinacota'
$ws.Range("P352").Value = 1083
$ws.Range("Q352").Value = 18
$ws.Range("R352").Value = 'Hortaliza'

# Row 353: Tomate, Larga vida, Segunda
$ws.Range("A353").Value = 11
$ws.Range("B353").Value = 'Vega Monumental Concepción'
$ws.Range("C353").Value = 'Bíobío'
$ws.Range("D353").Value = 44694
$ws.Range("E353").Value = 8
$ws.Range("F353").Value = 100112020
$ws.Range("G353").Value = 'Tomate'
$ws.Range("H353").Value = 'Larga vida'
$ws.Range("I353").Value = 'Segunda'
$ws.Range("J353").Value = 400
$ws.Range("K353").Value = 18000
$ws.Range("L353").Value = 18000
$ws.Range("M353").Value = 18000
$ws.Range("N353").Value = '$/bandeja 18 kilos'
$ws.Range("O353").Value = 'Región de Arica y Parinacota'
$ws.Range("P353").Value = 1000
$ws.Range("Q353").Value = 18
$ws.Range("R353").Value = 'Hortaliza'
